$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new trailing data row (row 12, "2021年") with the same
# look & feel as the existing rows, so first clone the formatting of the
# last existing row (row 11) down into row 12 ...
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)

# ... then fill in the new row's values.
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 118
$ws.Range("D12").Value = 910
$ws.Range("F12").Value = 32
$ws.Range("G12").Value = 760

# Columns C and E have no data for 2021 (matching the blank-but-present
# cells used throughout the rest of the table for missing values), so give
# them an empty value ...
$ws.Range("C12").Value = "'"
$ws.Range("E12").Value = "'"

# ... then restore the plain (non quote-prefixed) look of an ordinary blank
# data cell, copying it from an existing blank cell in the same column.
$ws.Range("C3").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E12").PasteSpecial(-4122)
